# This script applies the "deeper hypertuning" update to the results sheet.
# It updates the best_params / best_model strings, the confusion matrices,
# and the numeric score / confusion-matrix / training-time columns for the
# Random Forest, LightGBM and XGBoost rows (and the score/training-time for
# the CART row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: CART (best_params / best_model / confusion_matrix unchanged) ---
$ws.Range("B2").Value = "{'max_depth': 5, 'min_samples_leaf': 1, 'min_samples_split': 2}"
$ws.Range("C2").Value = 0.8474499348264011
$ws.Range("M2").Value = 17.09272265434265

# --- Row 3: Random Forest ---
$ws.Range("B3").Value = "{'max_depth': 25, 'min_samples_split': 2, 'n_estimators': 150}"
$ws.Range("C3").Value = 0.8694869060315205
$ws.Range("D3").Value = "RandomForestClassifier(max_depth=25, n_estimators=150)"
$ws.Range("E3").Value = "[[211  66]`n [ 23 426]]"
$ws.Range("F3").Value = 426
$ws.Range("G3").Value = 66
$ws.Range("H3").Value = 23
$ws.Range("I3").Value = 211
$ws.Range("J3").Value = 0.8795341555873707
$ws.Range("K3").Value = 0.8774104683195593
$ws.Range("L3").Value = 0.8750535215306668
$ws.Range("M3").Value = 694.8173654079437

# --- Row 4: LightGBM ---
$ws.Range("B4").Value = "{'learning_rate': 0.2, 'n_estimators': 50, 'num_leaves': 31}"
$ws.Range("C4").Value = 0.8746474700793933
$ws.Range("D4").Value = "LGBMClassifier(learning_rate=0.2, n_estimators=50)"
$ws.Range("E4").Value = "[[227  50]`n [ 33 416]]"
$ws.Range("F4").Value = 416
$ws.Range("H4").Value = 33
$ws.Range("J4").Value = 0.8852153471447268
$ws.Range("K4").Value = 0.8856749311294766
$ws.Range("L4").Value = 0.8849273526042128
$ws.Range("M4").Value = 237.3406167030334

# --- Row 5: XGBoost ---
$ws.Range("B5").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 100}"
$ws.Range("C5").Value = 0.8725844294347672
$ws.Range("D5").Value = "XGBClassifier(base_score=None, booster=None, callbacks=None,`n              colsample_bylevel=None, colsample_bynode=None,`n              colsample_bytree=None, device=None, early_stopping_rounds=None,`n              enable_categorical=True, eval_metric=None, feature_types=None,`n              gamma=None, grow_policy=None, importance_type=None,`n              interaction_constraints=None, learning_rate=0.1, max_bin=None,`n              max_cat_threshold=None, max_cat_to_onehot=None,`n              max_delta_step=None, max_depth=3, max_leaves=None,`n              min_child_weight=None, missing=nan, monotone_constraints=None,`n              multi_strategy=None, n_estimators=100, n_jobs=None,`n              num_parallel_tree=None, random_state=None, ...)"
$ws.Range("E5").Value = "[[219  58]`n [ 23 426]]"
$ws.Range("F5").Value = 426
$ws.Range("G5").Value = 58
$ws.Range("H5").Value = 23
$ws.Range("I5").Value = 219
$ws.Range("J5").Value = 0.8896250256130045
$ws.Range("K5").Value = 0.8884297520661157
$ws.Range("L5").Value = 0.8867604390243378
$ws.Range("M5").Value = 290.6994128227234

# Re-fit row heights for the rows whose multi-line cell text changed so the
# saved file doesn't pick up stray explicit row heights (matches the diff,
# which leaves the <row> elements untouched).
$ws.Rows(3).EntireRow.AutoFit()
$ws.Rows(4).EntireRow.AutoFit()
$ws.Rows(5).EntireRow.AutoFit()
